$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.633.19"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "2.236.23"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'305.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'94.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").Value = "'0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("D10").Value = "'34.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").Value = "'0.0803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("D12").Value = "'7.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").Value = "2.577.89"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").Value = "2.227.75"
$ws.Range("E15").Value = "  -4.09%  "

$ws.Range("D16").Value = "'0.833"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").Value = "'13.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "44.385.45"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("E19").Value = "  -3.11%  "

$ws.Range("D20").Value = "'11.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.88%  "

$ws.Range("D21").Value = "'6.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "

$ws.Range("D22").Value = "'65.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").Value = "'239.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("D24").Value = "'2.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.75%  "

$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +4.03%  "

$ws.Range("D28").Value = "'9.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.30%  "

$ws.Range("D29").Value = "'37.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").Value = "'5.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "

$ws.Range("E31").Value = "  -1.49%  "

$ws.Range("D32").Value = "'150.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("E33").Value = "  -1.78%  "

$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").Value = "'3.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.64%  "

$ws.Range("D36").Value = "'0.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("D38").Value = "'1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.50%  "

$ws.Range("D39").Value = "'15.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("E40").Value = "  -3.18%  "

$ws.Range("D41").Value = "'0.0301"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("E42").Value = "  -3.67%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "1.830.49"
$ws.Range("E44").Value = "  +5.21%  "

$ws.Range("D45").Value = "'1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.31%  "

$ws.Range("D46").Value = "'80.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.25%  "

$ws.Range("D47").Value = "'0.189"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").Value = "'98.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("E49").Value = "  -1.77%  "

$ws.Range("D50").Value = "'69.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").Value = "'7.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.93%  "
